# "update employer quick registration"
# Quick-registration sheet tweaks:
#  - shorten the brand address, append the "||," delimiter marker to the
#    uploaded-photo path
#  - widen the BrandLocation (B) and Media (G) columns so the longer values
#    are readable
#  - scroll the sheet over a bit and move the active selection from H3 to G4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates -------------------------------------------------
$ws.Range("B2").Value = "Ha-Banai St 29, Holon, 58857, Israel"
$ws.Range("G2").Value = "C:\Users\Harri\Desktop\Food.jpg||,"

# --- Column width updates --------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 42.85546875
$ws.Columns.Item(7).ColumnWidth = 53.140625

# --- View state: scroll the visible window and move the selection ---------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("G4").Select() | Out-Null
